$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their original Text storage type,
# so values like "22.30" or "6.00" are not silently coerced to numbers
# (which would drop the trailing zero / change the cell type).
$priceCells = @("D2", "D3", "D5", "D7", "D10", "D12", "D13", "D15", "D16", "D17", "D18", "D19", "D21", "D25", "D26", "D27", "D28", "D29", "D34", "D36", "D37", "D38", "D40", "D44", "D45", "D47", "D48", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "36.659.21"
$ws.Range("E2").Value = "  +0.38%  "

$ws.Range("D3").Value = "1.964.76"
$ws.Range("E3").Value = "  +1.03%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "244.56"
$ws.Range("E5").Value = "  +0.43%  "

$ws.Range("E6").Value = "  +0.98%  "

$ws.Range("D7").Value = "59.33"
$ws.Range("E7").Value = "  +1.59%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  +2.72%  "

$ws.Range("D10").Value = "0.0813"
$ws.Range("E10").Value = "  -2.75%  "

$ws.Range("E11").Value = "  -0.38%  "

$ws.Range("D12").Value = "22.30"
$ws.Range("E12").Value = "  +3.21%  "

$ws.Range("D13").Value = "2.253.74"
$ws.Range("E13").Value = "  +1.09%  "

$ws.Range("E14").Value = "  +1.11%  "

$ws.Range("D15").Value = "13.73"
$ws.Range("E15").Value = "  +1.41%  "

$ws.Range("D16").Value = "5.27"
$ws.Range("E16").Value = "  +0.44%  "

$ws.Range("D17").Value = "1.967.99"
$ws.Range("E17").Value = "  +0.73%  "

$ws.Range("D18").Value = "36.535.03"
$ws.Range("E18").Value = "  +0.37%  "

$ws.Range("D19").Value = "70.05"
$ws.Range("E19").Value = "  +0.48%  "

$ws.Range("E20").Value = "  -0.41%  "

$ws.Range("D21").Value = "229.27"
$ws.Range("E21").Value = "  -0.20%  "

$ws.Range("E22").Value = "  +0.22%  "

$ws.Range("E23").Value = "  -0.24%  "

$ws.Range("E24").Value = "  +0.78%  "

$ws.Range("D25").Value = "2.34"
$ws.Range("E25").Value = "  +2.54%  "

$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").Value = "0.141"
$ws.Range("E26").Value = "  +9.10%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "9.22"
$ws.Range("E27").Value = "  +0.06%  "

$ws.Range("D28").Value = "160.16"
$ws.Range("E28").Value = "  -0.97%  "

$ws.Range("D29").Value = "19.39"
$ws.Range("E29").Value = "  -0.19%  "

$ws.Range("E30").Value = "  +1.84%  "

$ws.Range("E31").Value = "  -0.02%  "

$ws.Range("E32").Value = "  +1.07%  "

$ws.Range("E33").Value = "  -1.10%  "

$ws.Range("D34").Value = "4.27"
$ws.Range("E34").Value = "  +0.51%  "

$ws.Range("E35").Value = "  +0.04%  "

$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "3.42"
$ws.Range("E36").Value = "  +12.91%  "

$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "2.25"
$ws.Range("E37").Value = "  +5.72%  "

$ws.Range("D38").Value = "6.00"
$ws.Range("E38").Value = "  -3.59%  "

$ws.Range("E39").Value = "  +0.22%  "

$ws.Range("D40").Value = "0.0986"
$ws.Range("E40").Value = "  +0.35%  "

$ws.Range("E41").Value = "  +1.87%  "

$ws.Range("E42").Value = "  -0.02%  "

$ws.Range("E43").Value = "  +1.00%  "

$ws.Range("D44").Value = "16.10"
$ws.Range("E44").Value = "  +0.32%  "

$ws.Range("D45").Value = "1.361.47"
$ws.Range("E45").Value = "  +0.89%  "

$ws.Range("E46").Value = "  +0.50%  "

$ws.Range("D47").Value = "87.88"
$ws.Range("E47").Value = "  +0.25%  "

$ws.Range("D48").Value = "7.13"
$ws.Range("E48").Value = "  +0.22%  "

$ws.Range("E49").Value = "  +0.57%  "

$ws.Range("D50").Value = "2.144.77"
$ws.Range("E50").Value = "  +1.14%  "

$ws.Range("D51").Value = "43.87"
$ws.Range("E51").Value = "  -3.32%  "
